$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # TRD
$ws2 = $wb.Worksheets.Item(2)   # Diseño
$ws3 = $wb.Worksheets.Item(3)   # Recursos

# ---------------------------------------------------------------------------
# TRD sheet (sheet1): widen column C, fill in the new task rows (4-24) and a
# SUBTOTAL totals row (25) driven off the table.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 24.7

$ws1.Range("A3").Value = 2

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = "Lista de tareas"
$ws1.Range("D4").Value = 0.5

$ws1.Range("A5").Value = 2
$ws1.Range("C5").Value = "Crear tarea"
$ws1.Range("D5").Value = 0.5

$ws1.Range("C6").Value = "Actualizar tarea"
$ws1.Range("D6").Value = 0.5

$ws1.Range("C7").Value = "Comentar"
$ws1.Range("D7").Value = 0.5

$ws1.Range("C8").Value = "Borrar"
$ws1.Range("D8").Value = 0.5

$ws1.Range("C9").Value = "Base de datos"
$ws1.Range("D9").Value = 0.5

$ws1.Range("C10").Value = "Tabla Tareas"
$ws1.Range("D10").Value = 0.5

$ws1.Range("C11").Value = "Tabla Comentarios"
$ws1.Range("D11").Value = 0.5

$ws1.Range("C12").Value = "Ltareas - Leer todos los datos de la base de datos"
$ws1.Range("D12").Value = 0.5

$ws1.Range("C13").Value = "genera html de lista de tareas"
$ws1.Range("D13").Value = 0.5

$ws1.Range("C14").Value = "Incluir menus de navegacion."
$ws1.Range("D14").Value = 0.5

$ws1.Range("C15").Value = "Formulario de Crear tarea (html)"
$ws1.Range("D15").Value = 0.5

$ws1.Range("C16").Value = "Guardar la tarea en base de datos"
$ws1.Range("D16").Value = 0.5

$ws1.Range("C17").Value = "Formulario de actualizar tarea (Html)"
$ws1.Range("D17").Value = 0.5

$ws1.Range("C18").Value = "Actualizar la tarea en db."
$ws1.Range("D18").Value = 0.5

$ws1.Range("C19").Value = "Formulatio borrar tarea (html)"
$ws1.Range("D19").Value = 0.5

$ws1.Range("C20").Value = "Borrar tarea de db."
$ws1.Range("D20").Value = 0.5

$ws1.Range("C21").Value = "Listar comentarios (html)"
$ws1.Range("D21").Value = 0.5

$ws1.Range("C22").Value = "formulario de Comentarios"
$ws1.Range("D22").Value = 0.5

$ws1.Range("C23").Value = "Agregar comentario en DB."
$ws1.Range("D23").Value = 0.5

$ws1.Range("C24").Value = "Diseñar interface (CSS). "
$ws1.Range("D24").Value = 0.5

# Grow the Table1 listobject to cover the new rows, then turn on a totals row
# with a SUM (SUBTOTAL) formula under "Tiempo Estimado".
$tbl = $ws1.ListObjects.Item(1)
$tbl.Resize($ws1.Range("A1:E24"))
$tbl.ShowTotals = $true
$col = $tbl.ListColumns.Item(4)
$col.TotalsCalculation = 1
$ws1.Range("D25").Formula = "=SUBTOTAL(109,Table1[Tiempo Estimado])"

# ---------------------------------------------------------------------------
# Recursos sheet (sheet3): the free-text hint cells become real numbers now
# that the estimate/example is filled in.
# ---------------------------------------------------------------------------
$ws3.Range("B3").Value = 11.5
$ws3.Range("B4").Value = 100
$ws3.Range("B5").Value = 0
$ws3.Range("B7").Value = 0
$ws3.Range("B8").Value = 66

$ws3.Columns.Item(2).ColumnWidth = 21.5

# ---------------------------------------------------------------------------
# Sheet views / selections / active sheet.
# ---------------------------------------------------------------------------
$ws2.Range("B8").Select()
$ws2.Range("D2").Select()

$ws1.Activate()
$excel.ActiveWindow.Zoom = 115
$ws1.Range("D29").Select()

$ws3.Activate()
$ws3.Range("B4").Select()
